$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as text (e.g. "174.40", "9.80", "2.00").
# Plain decimal-looking strings would otherwise be auto-coerced into numbers
# by Excel (dropping trailing zeros / changing "2.00" -> 2), so force those
# specific cells to remain text before assigning the new values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '66.438.64'
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("D3").Value = '3.440.05'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '580.74'
$ws.Range("E5").Value = '  -2.86%  '
$ws.Range("D6").Value = '174.40'
$ws.Range("E6").Value = '  -2.85%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").Value = '3.437.48'
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("E10").Value = '  -3.93%  '
$ws.Range("D11").Value = '6.86'
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("E12").Value = '  -3.05%  '
$ws.Range("D13").Value = '4.032.38'
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").Value = '30.92'
$ws.Range("E14").Value = '  -3.95%  '
$ws.Range("E15").Value = '  -3.59%  '
$ws.Range("D16").Value = '66.417.99'
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("E17").Value = '  -3.70%  '
$ws.Range("D18").Value = '3.438.01'
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("E19").Value = '  -4.68%  '
$ws.Range("E20").Value = '  -4.25%  '
$ws.Range("D21").Value = '373.51'
$ws.Range("E21").Value = '  -5.31%  '
$ws.Range("D22").Value = '7.74'
$ws.Range("E22").Value = '  -2.22%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = '70.76'
$ws.Range("E25").Value = '  -3.31%  '
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = '9.80'
$ws.Range("E28").Value = '  -6.66%  '
$ws.Range("E29").Value = '  -2.81%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").Value = '5.86'
$ws.Range("E31").Value = '  -5.05%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '23.81'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '2.00'
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("E34").Value = '  -6.74%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  -4.01%  '
$ws.Range("E37").Value = '  -5.80%  '
$ws.Range("D38").Value = '160.03'
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '27.05'
$ws.Range("E40").Value = '  +3.45%  '
$ws.Range("E41").Value = '  -5.66%  '
$ws.Range("E42").Value = '  -4.68%  '
$ws.Range("D43").Value = '6.51'
$ws.Range("E43").Value = '  -5.49%  '
$ws.Range("E44").Value = '  -5.15%  '
$ws.Range("D45").Value = '2.678.81'
$ws.Range("E45").Value = '  -6.06%  '
$ws.Range("D46").Value = '0.0692'
$ws.Range("E46").Value = '  -5.11%  '
$ws.Range("E47").Value = '  -5.36%  '
$ws.Range("D48").Value = '40.24'
$ws.Range("E48").Value = '  -4.43%  '
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("D50").Value = '318.84'
$ws.Range("E50").Value = '  -5.66%  '
$ws.Range("E51").Value = '  -5.06%  '
